# #CRM-41 Add acknowledge date in BB Adv Search form
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before B ("Tracking ID") and a new column before H
# ("Acknowledge Date"), shifting the existing columns to the right - this
# matches inserting the two new fields into the Buyback Order Snapshot
# advanced-search export.

# Insert "Tracking ID" right after "Order ID" (new column B)
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "Tracking ID"
$ws.Range("B2").Value = "{order:tracking_id}"

# Insert "Acknowledge Date" right after "Delivery Date" (new column H,
# since the previous insert shifted the original columns one to the right)
$ws.Columns("H").Insert()
$ws.Range("H1").Value = "Acknowledge Date"
$ws.Range("H2").Value = "{order:acknowledge_date}"
